# Swap the data (columns B through AD) between specific pairs of adjacent
# rows. Column A (the sequential row index) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(39, 40),
    @(60, 61),
    @(111, 112),
    @(132, 133),
    @(186, 187),
    @(201, 202),
    @(230, 231)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
